$d = $word.ActiveDocument

function New-PkgXml($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- 1. Title paragraph: "Neon" + "Cluster" -> proofErr-wrapped "neonCLUSTER" ---
$titleBody = '<w:p><w:pPr><w:pStyle w:val="Title"/><w:jc w:val="center"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>neonCLUSTER</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Depl</w:t></w:r><w:r><w:t>oyment</w:t></w:r></w:p>'
$d.Paragraphs(1).Range.InsertXML((New-PkgXml $titleBody))

# --- 2. "...to create NeonCluster development and production clusters..." ---
$overviewBody = '<w:p><w:r><w:t xml:space="preserve">This document </w:t></w:r><w:r><w:t>starts out by describing</w:t></w:r><w:r><w:t xml:space="preserve"> how the scripts and other assets located within this source folder are to be used to </w:t></w:r><w:r><w:t>create</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neonCLUSTER</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>development and production</w:t></w:r><w:r><w:t xml:space="preserve"> cluster</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space="preserve"> on Ubuntu 1</w:t></w:r><w:r><w:t>6</w:t></w:r><w:r><w:t>.04 LTS.</w:t></w:r></w:p>'
$d.Paragraphs(6).Range.InsertXML((New-PkgXml $overviewBody))

# --- 3. "Creating a Development NeonCluster" heading ---
$headingBody = '<w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t xml:space="preserve">Creating a </w:t></w:r><w:r><w:t xml:space="preserve">Development </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>neonCLUSTER</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$d.Paragraphs(8).Range.InsertXML((New-PkgXml $headingBody))

# --- 4. Move the _GoBack bookmark from next to "neon-cli" up to the empty
#        paragraph pair near the top of the document. ---
$d.Bookmarks("_GoBack").Delete()
$d.Bookmarks.Add("_GoBack", $d.Paragraphs(4).Range)
